$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("A3").Value = "WH/IN/3"
$ws.Range("B3").Value = "Camptocamp"
$ws.Range("C3").Value = "ABCD"
$ws.Range("D3").Value = 42952
$ws.Range("D3").NumberFormat = "DD/MM/YY"
$ws.Range("E3").Value = "CONS_DEL03"
$ws.Range("F3").Value = 20

# --- Row 4 ---
$ws.Range("A4").Value = "WH/IN/4"
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("B4").Value = "Think Big Systems"
$ws.Range("C4").Value = "SO006"
$ws.Range("D4").Value = 36746
$ws.Range("D4").NumberFormat = "DD/MM/YY"
$ws.Range("E4").Value = "E-COM01"
$ws.Range("F4").Value = 2

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 16.569727891156468
$ws.Columns.Item(5).ColumnWidth = 14.447278911564666

# --- Selection ---
[void]$ws.Range("D8").Select()

Write-Output "done"
